$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 18 (impediment / LOSE) - this shifts rows 19,20 up to 18,19
# and removes the now-unused "impediment" shared string.
$ws.Rows.Item(18).Delete()

# The former row 19 (resources/LOSE) is now row 18; restyle its B cell
# to match the diff (s=3 -> automatic/no-theme font color).
$ws.Cells.Item(18, 2).Font.ColorIndex = -4105

# Update the selected cell shown in the sheet view.
$ws.Range("B21").Select()

# Ensure page setup matches (portrait orientation).
$ws.PageSetup.Orientation = 1
